$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")
$ws.Range("Z100").Font.ThemeColor = 5
$ws.Range("Z100").Font.TintAndShade = [double]-0.249977111117893
Write-Host ($ws.Range("Z100").Font.TintAndShade)
